$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024751861797291
$ws.Range("D2").Value = 1.047194623822245
$ws.Range("E2").Value = 1.036794980570205
$ws.Range("F2").Value = 1.051681288936122
$ws.Range("I2").Value = 1.040296443461219
$ws.Range("J2").Value = 1.029924720049001
$ws.Range("K2").Value = 1.049957942301966
$ws.Range("L2").Value = 1.039587701869937
$ws.Range("M2").Value = 1.054432116978946
$ws.Range("N2").Value = 1.014009659332276
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025644096042755
$ws.Range("D3").Value = 1.047768494738549
$ws.Range("E3").Value = 1.037551551465974
$ws.Range("F3").Value = 1.052437306499862
$ws.Range("I3").Value = 1.040417897962093
$ws.Range("J3").Value = 1.030456290935767
$ws.Range("K3").Value = 1.050344225848831
$ws.Range("L3").Value = 1.040154088676741
$ws.Range("M3").Value = 1.055000975048345
$ws.Range("N3").Value = 1.014186939071788
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026222052866931
$ws.Range("D4").Value = 1.048138796009451
$ws.Range("E4").Value = 1.03804164884026
$ws.Range("F4").Value = 1.052926190901132
$ws.Range("I4").Value = 1.04049432402437
$ws.Range("J4").Value = 1.030800254560706
$ws.Range("K4").Value = 1.050592445620486
$ws.Range("L4").Value = 1.040520461358431
$ws.Range("M4").Value = 1.055368091356972
$ws.Range("N4").Value = 1.014301606395667
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026465173370303
$ws.Range("D5").Value = 1.048294221426423
$ws.Range("E5").Value = 1.038247814534767
$ws.Range("F5").Value = 1.053131641598045
$ws.Range("I5").Value = 1.040525934515635
$ws.Range("J5").Value = 1.030944856095377
$ws.Range("K5").Value = 1.050696380824806
$ws.Range("L5").Value = 1.040674454851715
$ws.Range("M5").Value = 1.055522192216235
$ws.Range("N5").Value = 1.014349801447574
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026506002979888
$ws.Range("D6").Value = 1.04832030335594
$ws.Range("E6").Value = 1.038282438159575
$ws.Range("F6").Value = 1.053166133135045
$ws.Range("I6").Value = 1.040531211586067
$ws.Range("J6").Value = 1.030969135249476
$ws.Range("K6").Value = 1.050713807531343
$ws.Range("L6").Value = 1.040700309259454
$ws.Range("M6").Value = 1.055548052601453
$ws.Range("N6").Value = 1.014357892945872
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026225300877308
$ws.Range("D7").Value = 1.048140873793719
$ws.Range("E7").Value = 1.038044403132729
$ws.Range("F7").Value = 1.052928936446039
$ws.Range("I7").Value = 1.040494748446213
$ws.Range("J7").Value = 1.030802186738025
$ws.Range("K7").Value = 1.050593836046274
$ws.Range("L7").Value = 1.04052251914512
$ws.Range("M7").Value = 1.055370151385397
$ws.Range("N7").Value = 1.014302250424436
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025053267013553
$ws.Range("D8").Value = 1.047388778420787
$ws.Range("E8").Value = 1.037050553226229
$ws.Range("F8").Value = 1.051936851376486
$ws.Range("I8").Value = 1.040337936767355
$ws.Range("J8").Value = 1.03010436585302
$ws.Range("K8").Value = 1.050088846043516
$ws.Range("L8").Value = 1.039779138385611
$ws.Range("M8").Value = 1.054624565520817
$ws.Range("N8").Value = 1.014069580809314
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022992814265424
$ws.Range("D9").Value = 1.046055689150575
$ws.Range("E9").Value = 1.035303519777265
$ws.Range("F9").Value = 1.050186397028776
$ws.Range("I9").Value = 1.040045098722995
$ws.Range("J9").Value = 1.028874781533782
$ws.Range("K9").Value = 1.049185809170243
$ws.Range("L9").Value = 1.038468371238158
$ws.Range("M9").Value = 1.05330338098538
$ws.Range("N9").Value = 1.013659265264153
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021622500048253
$ws.Range("D10").Value = 1.045161851142519
$ws.Range("E10").Value = 1.034141800057886
$ws.Range("F10").Value = 1.049018028272015
$ws.Range("I10").Value = 1.039838837123775
$ws.Range("J10").Value = 1.028055172168267
$ws.Range("K10").Value = 1.048575031729833
$ws.Range("L10").Value = 1.037594050925655
$ws.Range("M10").Value = 1.052417750776899
$ws.Range("N10").Value = 1.013385530576515
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021029943072578
$ws.Range("D11").Value = 1.044773624755862
$ws.Range("E11").Value = 1.033639489260708
$ws.Range("F11").Value = 1.048511806851644
$ws.Range("I11").Value = 1.039746920029551
$ws.Range("J11").Value = 1.027700313699183
$ws.Range("K11").Value = 1.048308504878844
$ws.Range("L11").Value = 1.037215365737686
$ws.Range("M11").Value = 1.052033138468846
$ws.Range("N11").Value = 1.013266960451307
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020809962235665
$ws.Range("D12").Value = 1.044629243822153
$ws.Range("E12").Value = 1.033453018858629
$ws.Range("F12").Value = 1.048323729247112
$ws.Range("I12").Value = 1.039712387731511
$ws.Range("J12").Value = 1.027568510291721
$ws.Range("K12").Value = 1.048209197799105
$ws.Range("L12").Value = 1.037074691793066
$ws.Range("M12").Value = 1.051890108631809
$ws.Range("N12").Value = 1.013222912459851
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020857143366183
$ws.Range("D13").Value = 1.044660221968273
$ws.Range("E13").Value = 1.03349301237764
$ws.Range("F13").Value = 1.048364074507877
$ws.Range("I13").Value = 1.039719812666432
$ws.Range("J13").Value = 1.027596782242634
$ws.Range("K13").Value = 1.048230513375165
$ws.Range("L13").Value = 1.037104867415469
$ws.Range("M13").Value = 1.05192079657654
$ws.Range("N13").Value = 1.013232361157778
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021011756907869
$ws.Range("D14").Value = 1.04476169376798
$ws.Range("E14").Value = 1.033624073300838
$ws.Range("F14").Value = 1.048496261192146
$ws.Range("I14").Value = 1.039744073532628
$ws.Range("J14").Value = 1.02768941864248
$ws.Range("K14").Value = 1.048300302379993
$ws.Range("L14").Value = 1.037203737850767
$ws.Range("M14").Value = 1.052021318988263
$ws.Range("N14").Value = 1.013263319545814
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021107035513804
$ws.Range("D15").Value = 1.044824190588095
$ws.Range("E15").Value = 1.033704838916426
$ws.Range("F15").Value = 1.048577699947769
$ws.Range("I15").Value = 1.039758969780694
$ws.Range("J15").Value = 1.02774649592492
$ws.Range("K15").Value = 1.048343261043061
$ws.Range("L15").Value = 1.037264653455387
$ws.Range("M15").Value = 1.052083231997058
$ws.Range("N15").Value = 1.013282393278774
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021661843029556
$ws.Range("D16").Value = 1.045187591595241
$ws.Range("E16").Value = 1.034175152125909
$ws.Range("F16").Value = 1.049051618172199
$ws.Range("I16").Value = 1.039844882579255
$ws.Range("J16").Value = 1.028078723850685
$ws.Range("K16").Value = 1.048592677069847
$ws.Range("L16").Value = 1.037619181056408
$ws.Range("M16").Value = 1.052443252609866
$ws.Range("N16").Value = 1.013393398853268
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022010073723221
$ws.Range("D17").Value = 1.045415226782425
$ws.Range("E17").Value = 1.034470361614441
$ws.Range("F17").Value = 1.049348813070028
$ws.Range("I17").Value = 1.039898076895942
$ws.Range("J17").Value = 1.028287132582827
$ws.Range("K17").Value = 1.048748579915069
$ws.Range("L17").Value = 1.037841541432783
$ws.Range("M17").Value = 1.05266878318934
$ws.Range("N17").Value = 1.013463018985682
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022213267715461
$ws.Range("D18").Value = 1.045547887703018
$ws.Range("E18").Value = 1.034642621809614
$ws.Range("F18").Value = 1.049522131768135
$ws.Range("I18").Value = 1.039928852875454
$ws.Range("J18").Value = 1.028408697440041
$ws.Range("K18").Value = 1.048839316896248
$ws.Range("L18").Value = 1.03797123085761
$ws.Range("M18").Value = 1.052800222345003
$ws.Range("N18").Value = 1.013503623227701
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022282564623314
$ws.Range("D19").Value = 1.045593102103509
$ws.Range("E19").Value = 1.034701369817409
$ws.Range("F19").Value = 1.049581223801307
$ws.Range("I19").Value = 1.039939304019753
$ws.Range("J19").Value = 1.02845014848529
$ws.Range("K19").Value = 1.048870222136738
$ws.Range("L19").Value = 1.03801544992465
$ws.Range("M19").Value = 1.052845021152784
$ws.Range("N19").Value = 1.013517467533024
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021972703899785
$ws.Range("D20").Value = 1.045390815542826
$ws.Range("E20").Value = 1.034438681228963
$ws.Range("F20").Value = 1.049316929963534
$ws.Range("I20").Value = 1.039892395641145
$ws.Range("J20").Value = 1.028264771920665
$ws.Range("K20").Value = 1.048731873535215
$ws.Range("L20").Value = 1.037817685253568
$ws.Range("M20").Value = 1.052644597155104
$ws.Range("N20").Value = 1.013455549815381
$ws.Range("B21").Value = 1.019999999999999
$ws.Range("C21").Value = 1.020966223717136
$ws.Range("D21").Value = 1.044731817703045
$ws.Range("E21").Value = 1.033585476069311
$ws.Range("F21").Value = 1.048457336717568
$ws.Range("I21").Value = 1.039736940074825
$ws.Range("J21").Value = 1.027662139333925
$ws.Range("K21").Value = 1.048279759721886
$ws.Range("L21").Value = 1.037174623330464
$ws.Range("M21").Value = 1.051991722257995
$ws.Range("N21").Value = 1.013254203227987
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020334111591194
$ws.Range("D22").Value = 1.044316460275621
$ws.Range("E22").Value = 1.033049670549707
$ws.Range("F22").Value = 1.047916620325536
$ws.Range("I22").Value = 1.039636941966399
$ws.Range("J22").Value = 1.027283280379571
$ws.Range("K22").Value = 1.047993721673283
$ws.Range("L22").Value = 1.036770227921581
$ws.Range("M22").Value = 1.051580264540616
$ws.Range("N22").Value = 1.013127575475306
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020669139222316
$ws.Range("D23").Value = 1.044536744994369
$ws.Range("E23").Value = 1.033333650080498
$ws.Range("F23").Value = 1.048203287867311
$ws.Range("I23").Value = 1.039690166433096
$ws.Range("J23").Value = 1.027484116498282
$ws.Range("K23").Value = 1.048145523623756
$ws.Range("L23").Value = 1.036984612389439
$ws.Range("M23").Value = 1.051798477329307
$ws.Range("N23").Value = 1.013194706245671
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021989589476875
$ws.Range("D24").Value = 1.045401846286275
$ws.Range("E24").Value = 1.034452996016962
$ws.Range("F24").Value = 1.049331336660058
$ws.Range("I24").Value = 1.039894963532415
$ws.Range("J24").Value = 1.028274875730477
$ws.Range("K24").Value = 1.048739423042688
$ws.Range("L24").Value = 1.037828464863949
$ws.Range("M24").Value = 1.05265552611992
$ws.Range("N24").Value = 1.013458924824347
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023524911191124
$ws.Range("D25").Value = 1.046401235354411
$ws.Range("E25").Value = 1.035754654256371
$ws.Range("F25").Value = 1.050639187479803
$ws.Range("I25").Value = 1.040122754528336
$ws.Range("J25").Value = 1.029192643629702
$ws.Range("K25").Value = 1.049420816484417
$ws.Range("L25").Value = 1.03880732605558
$ws.Range("M25").Value = 1.053645800223409
$ws.Range("N25").Value = 1.013765377190566
